$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "General" vs "Property" search scenarios (rows 32 & 35) ---
# The Given/When (columns A & B) content for these two rows was swapped while
# the Then (column C) result stayed tied to its row.
$a32 = $ws.Range("A32").Formula
$b32 = $ws.Range("B32").Formula
$a35 = $ws.Range("A35").Formula
$b35 = $ws.Range("B35").Formula

$ws.Range("A32").Formula = $a35
$ws.Range("B32").Formula = $b35
$ws.Range("A35").Formula = $a32
$ws.Range("B35").Formula = $b32

# --- Clean up flaky tests: mark rows 28-37 and row 60 as done, matching the
#     highlight formatting already used from row 38 onward. Copy the format
#     from an already-highlighted cell so the existing style is reused
#     instead of creating a duplicate style entry. ---
$ws.Range("A38:C38").Copy() | Out-Null
$ws.Range("A28:C37").PasteSpecial(-4122) | Out-Null
$ws.Range("A60:C60").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update the selection to reflect the last edited row ---
$ws.Range("A60:C60").Select()
